$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets 1-4 ("Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)") all share the
# same row layout (Fonte/Tecnologia header in col A, years in row 1).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Add the missing header in A1, reusing the existing header style (same
    # as B1) so no new style entry gets created.
    $ws.Range("B1").Copy($ws.Range("A1"))
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # Remove the bold/boxed style from the row labels in column A (A2:A12)
    # and fix the accented labels.
    $ws.Range("A2:A12").Style = "Normal"

    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."
}

# ---------------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)

$ws.Range("B1").Copy($ws.Range("A1"))
$ws.Range("A1").Value = "Período"

$ws.Range("A2:A3").Style = "Normal"
$ws.Range("A2").Value = "P.Médio"
$ws.Range("A3").Value = "P.Crítico"

# Remove the "Teto" row entirely.
$ws.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)

$ws.Range("B1").Copy($ws.Range("A1"))
$ws.Range("A1").Value = "Tipo Expansão"

# B1 must hold the text "2015" (same as the header used on the other
# sheets), not a number - copy it from a sheet that already has it so the
# correct string cell (and style) is produced without adding a new one.
$wsYears = $wb.Worksheets.Item(1)
$wsYears.Range("B1").Copy($ws.Range("B1"))

$ws.Range("A2:A3").Style = "Normal"
$ws.Range("A2").Value = "Expansão Centralizada"
$ws.Range("B2").Value = 692
$ws.Range("A3").Value = "Expansão por GD"
$ws.Range("B3").Value = 99
